# Support formula evaluation after inserting the data.
# Adds a new "Sheet1" after "data-source" that pulls a value from the
# data-source sheet via formula, and makes the new sheet the active tab.

$wb = $excel.ActiveWorkbook

# The existing sheet holding the raw/inserted data.
$dataSourceSheet = $wb.Worksheets.Item("data-source")

# Insert the new sheet right after "data-source" (keeps it first, new one second).
$newSheet = $wb.Worksheets.Add($null, $dataSourceSheet)

# Pull a value from the data-source sheet so it recalculates once data lands there.
$newSheet.Range("A1").Formula = "='data-source'!C2"

# Land the cursor on A2 (below the formula); selecting a cell on the new sheet
# also makes it the active sheet/tab.
$newSheet.Range("A2").Select() | Out-Null
